$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "30.274.45"
Set-TextValue "E2" "  -0.15%  "

# Row 3
Set-TextValue "D3" "1.915.07"
Set-TextValue "E3" "  -0.83%  "

# Row 4
Set-TextValue "D4" "0.9999"
Set-TextValue "E4" "  +0.08%  "

# Row 5
Set-TextValue "D5" "0.7425"
Set-TextValue "E5" "  -1.18%  "

# Row 6
Set-TextValue "D6" "243.96"
Set-TextValue "E6" "  -1.89%  "

# Row 7
Set-TextValue "D7" "0.9998"
Set-TextValue "E7" "  -0.01%  "

# Row 8
Set-TextValue "D8" "0.3154"
Set-TextValue "E8" "  -1.77%  "

# Row 9
Set-TextValue "D9" "27.26"
Set-TextValue "E9" "  -4.42%  "

# Row 10
Set-TextValue "D10" "0.07017"
Set-TextValue "E10" "  -1.19%  "

# Row 11
Set-TextValue "D11" "0.7841"
Set-TextValue "E11" "  -0.51%  "

# Row 12
Set-TextValue "D12" "0.07973"
Set-TextValue "E12" "  -0.34%  "

# Row 13
Set-TextValue "D13" "1.917.46"
Set-TextValue "E13" "  -0.79%  "

# Row 14
Set-TextValue "D14" "5.301"
Set-TextValue "E14" "  -1.44%  "

# Row 15
Set-TextValue "D15" "92.06"
Set-TextValue "E15" "  -2.72%  "

# Row 16
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D16" "14.39"
Set-TextValue "E16" "  -1.88%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D17" "30.258.42"
Set-TextValue "E17" "  -0.24%  "

# Row 18
Set-TextValue "D18" "246.62"
Set-TextValue "E18" "  -2.93%  "

# Row 19
Set-TextValue "D19" "5.859"
Set-TextValue "E19" "  +0.75%  "

# Row 20
Set-TextValue "D20" "0.000007859"
Set-TextValue "E20" "  -2.13%  "

# Row 21
Set-TextValue "D21" "2.175.72"
Set-TextValue "E21" "  -0.36%  "

# Row 22
Set-TextValue "D22" "0.9993"
Set-TextValue "E22" "  -0.04%  "

# Row 23
Set-TextValue "D23" "0.9996"
Set-TextValue "E23" "  -0.01%  "

# Row 24
Set-TextValue "D24" "6.677"
Set-TextValue "E24" "  -2.02%  "

# Row 25
Set-TextValue "D25" "9.478"
Set-TextValue "E25" "  -1.04%  "

# Row 26
Set-TextValue "D26" "165.42"
Set-TextValue "E26" "  +0.53%  "

# Row 27
Set-TextValue "D27" "19.09"
Set-TextValue "E27" "  -0.06%  "

# Row 28
Set-TextValue "D28" "0.1279"
Set-TextValue "E28" "  -4.29%  "

# Row 29
Set-TextValue "D29" "2.137"
Set-TextValue "E29" "  -8.31%  "

# Row 30
Set-TextValue "D30" "1.352"
Set-TextValue "E30" "  -0.71%  "

# Row 31
Set-TextValue "E31" "  +1.10%  "

# Row 32
Set-TextValue "D32" "4.342"
Set-TextValue "E32" "  -2.36%  "

# Row 33
Set-TextValue "D33" "4.097"
Set-TextValue "E33" "  -1.19%  "

# Row 34
Set-TextValue "D34" "0.05251"
Set-TextValue "E34" "  +2.02%  "

# Row 35
Set-TextValue "D35" "1.310"
Set-TextValue "E35" "  +1.98%  "

# Row 36
Set-TextValue "D36" "0.7546"
Set-TextValue "E36" "  +0.40%  "

# Row 37
Set-TextValue "D37" "2.760"
Set-TextValue "E37" "  -0.42%  "

# Row 38
Set-TextValue "D38" "0.01949"
Set-TextValue "E38" "  -0.99%  "

# Row 39
Set-TextValue "D39" "2.796"
Set-TextValue "E39" "  -0.14%  "

# Row 40
Set-TextValue "D40" "6.419"
Set-TextValue "E40" "  -0.01%  "

# Row 41
Set-TextValue "D41" "76.32"
Set-TextValue "E41" "  -2.47%  "

# Row 42
Set-TextValue "D42" "0.4523"
Set-TextValue "E42" "  +0.11%  "

# Row 43
Set-TextValue "D43" "1.962"
Set-TextValue "E43" "  -1.69%  "

# Row 44
Set-TextValue "D44" "0.9987"
Set-TextValue "E44" "  -0.12%  "

# Row 45
Set-TextValue "D45" "7.777"
Set-TextValue "E45" "  +2.78%  "

# Row 46
Set-TextValue "D46" "0.8331"
Set-TextValue "E46" "  -0.41%  "

# Row 47
Set-TextValue "D47" "101.47"
Set-TextValue "E47" "  -0.92%  "

# Row 48
Set-TextValue "D48" "9.905"
Set-TextValue "E48" "  +0.79%  "

# Row 49
Set-TextValue "D49" "2.110.96"
Set-TextValue "E49" "  +1.06%  "

# Row 50
Set-TextValue "D50" "37.24"
Set-TextValue "E50" "  -0.53%  "

# Row 51
Set-TextValue "D51" "0.1220"
Set-TextValue "E51" "  +1.31%  "

